$d = $word.ActiveDocument

# 1. Funding support sentence (contains old course name + framework code, replace whole sentence first
#    so it doesn't get disturbed by the standalone title replacement below)
$d.Content.Find.Execute(
    "We are applying for WSQ funding support for this new course PWM-Security: Security Risk Analysis (Assess and Address Security Risks) according to Security Risk Analysis Assess and Address Security Risks SEC-SRM-3002-1.1 under Security Framework.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We are applying for WSQ funding support for this new course Storytelling and Storyboarding with Generative AI according to AI Content Generation for Script Development MED-MED-3004-1.1 under Media Framework.", 2)

# 2. Title heading
$d.Content.Find.Execute(
    "PWM-Security: Security Risk Analysis (Assess and Address Security Risks)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Storytelling and Storyboarding with Generative AI", 2)

# 3. Performance gaps paragraph
$d.Content.Find.Execute(
    "A notable deficiency lies in the ability to thoroughly analyze security situations and apply prior experiences to new, evolving threats. Many professionals find it challenging to connect seemingly unrelated events to identify underlying risks. The ability to think critically and adapt to novel situations is crucial. Moreover, collaboration between security teams is often fragmented.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The industry struggles with effectively utilizing AI tools to enhance the quality and consistency of video scripts. Many creators find it challenging to refine AI-generated content to meet the desired tone and clarity, resulting in scripts that lack coherence and engagement. This gap hinders the ability to produce high-quality video content that aligns with audience expectations.", 2)

# 4. Why this course addresses training needs paragraph
$d.Content.Find.Execute(
    "This course promotes a holistic approach to security by teaching individuals how to relate current situations to past experiences, thereby enhancing their ability to identify subtle patterns and risks. It also fosters collaboration among security teams, ensuring comprehensive risk analysis, and facilitates effective communication, which ensures that different perspectives and expertise contribute to a more robust security framework. This can allow the teams to address complex risks and implement preventive strategies.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Through this course, participants gain the ability to refine and enhance AI-generated video scripts, improving clarity and narrative consistency. The training focuses on utilizing AI tools to achieve the desired tone and engagement, enabling creators to produce high-quality video content that resonates with audiences. This skill set is crucial for meeting industry standards and audience expectations.", 2)

# 5. Date field
$d.Content.Find.Execute(
    ": 06 March 2025",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": 08 May 2025", 2)
